# Insert a new weekly price record at row 277 of Sheet1, pushing the
# existing rows 277-361 down to 278-362 (dimension grows from R361 to R362).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(277).Insert()

$ws.Range("A277").Value = 8
$ws.Range("B277").Value = "Terminal La Palmera de La Serena"
$ws.Range("C277").Value = "Coquimbo"
$ws.Range("D277").Value = 44463
$ws.Range("E277").Value = 4
$ws.Range("F277").Value = 100112024
$ws.Range("G277").Value = "Choclo"
$ws.Range("H277").Value = "Dulce o Americano"
$ws.Range("I277").Value = "Primera"
$ws.Range("J277").Value = 600
$ws.Range("K277").Value = 35000
$ws.Range("L277").Value = 37000
$ws.Range("M277").Value = 36000
$ws.Range("N277").Value = "$/malla 70 unidades"
$ws.Range("O277").Value = "Región de Arica y Parinacota"
$ws.Range("P277").Value = 514
$ws.Range("Q277").Value = 70
$ws.Range("R277").Value = "Hortaliza"
